$d = $word.ActiveDocument

# =====================================================================
# This script applies, in bottom-to-top paragraph order (so earlier
# edits never shift the paragraph index of work still to be done):
#
#   1. End of document: two new sentences appended to the "AVL seemed
#      to be..." paragraph (before the _GoBack bookmark), plus a brand
#      new trailing paragraph ("Overall, it seems...").
#   2. Two new blank paragraphs + a centered "Analysis" heading
#      inserted right before the "AVL seemed to be..." paragraph.
#   3. Table cell "Range Query (n=1000)" / "Range Query (n=100)" text
#      normalised (re-typed so the runs collapse the way Word does).
#   4. Table cell (Delete / 2-5 Tree Time) "N/A" -> ".04953 seconds"
#   5. Table cell (Delete / AVL Tree Time) ".000281 seconds" re-typed.
#   6. Table cell (Search / 2-5 Tree Time) "N/A" -> "0.00028" on one
#      paragraph and "seconds" on a new paragraph beneath it.
#   7. Table cell (Search / AVL Tree Time) ".000199 seconds" re-typed.
# =====================================================================

# ---------------------------------------------------------------------
# 1. End-of-document narrative additions
# ---------------------------------------------------------------------
$pAvl = $d.Paragraphs.Item(83)

# Append the two new sentences immediately before the (hidden) _GoBack
# bookmark that trails the paragraph, so they land before it just like
# in the target revision.
$insPos = $pAvl.Range.End - 1
$ins1 = $d.Range($insPos, $insPos)
$ins1.InsertBefore(" The time for AVL and 2-5 Tree to perform search, insert, sort, and range search seem to be about the same.")

$pAvl2 = $d.Paragraphs.Item(83)
$insPos2 = $pAvl2.Range.End - 1
$ins2 = $d.Range($insPos2, $insPos2)
$ins2.InsertBefore(" The delete operation seemed for 2-5 took much longer time. This may be due to how delete was handled or because of which words were selected for testing.")

# New trailing paragraph after the "AVL seemed..." paragraph.
$pAvl3 = $d.Paragraphs.Item(83)
$pAvl3.Range.InsertParagraphAfter()
$pOverall = $d.Paragraphs.Item(84)
$pOverall.Range.Text = "Overall, it seems that BST has the quickest operations (probably because of the lack of the balancing factor AVL has). However, since this is processing only 100 words for each operation, the run time for a large sample size should result in similar times for BST, AVL, and 2-5. "

# ---------------------------------------------------------------------
# 2. New "Analysis" heading (centered, 12pt) with two blank lines above
# ---------------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(82)
$pBlank.Range.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs.Item(83)
$pBlank2.Range.InsertParagraphAfter()
$pBlank3 = $d.Paragraphs.Item(84)
$pBlank3.Range.InsertParagraphAfter()

$pAnalysis = $d.Paragraphs.Item(85)
$pAnalysis.Range.Text = "Analysis"

$pAnalysis2 = $d.Paragraphs.Item(85)
$pAnalysis2.Alignment = 1
$pAnalysis3 = $d.Paragraphs.Item(85)
$pAnalysis3.Range.Font.Size = 12
$pAnalysis4 = $d.Paragraphs.Item(85)
$pAnalysis4.Range.Font.SizeBi = 12

# ---------------------------------------------------------------------
# 3. "Range Query (n=1000)" / "Range Query (n=100)" run normalisation
# ---------------------------------------------------------------------
$pRq1000 = $d.Paragraphs.Item(72)
$pRq1000.Range.Text = "TEMP_RQ_1000"
$pRq1000b = $d.Paragraphs.Item(72)
$pRq1000b.Range.Text = "Range Query (n=1000)"

$pRq100 = $d.Paragraphs.Item(62)
$pRq100.Range.Text = "TEMP_RQ_100"
$pRq100b = $d.Paragraphs.Item(62)
$pRq100b.Range.Text = "Range Query (n=100)"

# ---------------------------------------------------------------------
# 4. Delete / 2-5 Tree Time: "N/A" -> ".04953 seconds"
# ---------------------------------------------------------------------
$p39 = $d.Paragraphs.Item(39)
$p39.Range.Text = ".04953 seconds"
$p39b = $d.Paragraphs.Item(39)
$p39b.Range.Font.Name = "Times New Roman"
$p39c = $d.Paragraphs.Item(39)
$p39c.Range.Font.SizeBi = 10
$p39d = $d.Paragraphs.Item(39)
$p39d.Range.Font.Size = 10

# ---------------------------------------------------------------------
# 5. Delete / AVL Tree Time: re-type ".000281 seconds" (run merge only)
# ---------------------------------------------------------------------
$p37 = $d.Paragraphs.Item(37)
$p37start = $p37.Range.Start
$p37sub = $d.Range($p37start + 2, $p37start + 7)
$p37sub.Text = "TEMP5"
$p37sub2 = $d.Range($p37start + 2, $p37start + 7)
$p37sub2.Text = "00281"

# ---------------------------------------------------------------------
# 6. Search / 2-5 Tree Time: "N/A" -> "0.00028" + new paragraph "seconds"
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "0.00028"

$p18b = $d.Paragraphs.Item(18)
$p18b.Range.Font.Name = "Times New Roman"
$p18c = $d.Paragraphs.Item(18)
$p18c.Range.Font.SizeBi = 10
$p18d = $d.Paragraphs.Item(18)
$p18d.Range.Font.Size = 10

$p18e = $d.Paragraphs.Item(18)
$p18e.Range.InsertParagraphAfter()

$p19 = $d.Paragraphs.Item(19)
$p19.Range.Text = "seconds"
$p19b = $d.Paragraphs.Item(19)
$p19b.Range.Font.Name = "Times New Roman"
$p19c = $d.Paragraphs.Item(19)
$p19c.Range.Font.SizeBi = 10
$p19d = $d.Paragraphs.Item(19)
$p19d.Range.Font.Size = 10

# ---------------------------------------------------------------------
# 7. Search / AVL Tree Time: re-type ".000199 seconds" (run merge only)
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$p16start = $p16.Range.Start
$p16sub = $d.Range($p16start + 2, $p16start + 8)
$p16sub.Text = "TEMP6Z"
$p16sub2 = $d.Range($p16start + 2, $p16start + 8)
$p16sub2.Text = "00199 "

Write-Output "Done."
